# Change data source to CDC
# - Fix F661 value (30 -> 0.3)
# - Truncate several percentage values in rows 664-672 to 3 decimals
# - Append new rows 673-681 with data for 2021-03-02 (serial 44257)
# - Extend dimension / defined name range from H672 to H681

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix existing values ---
$ws.Cells.Item(661, 6).Value = 0.3

$ws.Cells.Item(664, 6).Value = 0.562
$ws.Cells.Item(664, 8).Value = 0.548

$ws.Cells.Item(665, 6).Value = 6.826
$ws.Cells.Item(665, 8).Value = 8.191

$ws.Cells.Item(666, 6).Value = 8.948
$ws.Cells.Item(666, 8).Value = 10.81

$ws.Cells.Item(667, 6).Value = 9.908
$ws.Cells.Item(667, 8).Value = 11.55

$ws.Cells.Item(668, 6).Value = 11.05
$ws.Cells.Item(668, 8).Value = 12.51

$ws.Cells.Item(669, 6).Value = 18.49
$ws.Cells.Item(669, 8).Value = 10.93

$ws.Cells.Item(670, 6).Value = 29.74
$ws.Cells.Item(670, 8).Value = 26.75

$ws.Cells.Item(671, 6).Value = 14.44
$ws.Cells.Item(671, 8).Value = 18.69

$ws.Cells.Item(672, 6).Value = 0.012
$ws.Cells.Item(672, 8).Value = 0.002

# --- Append new rows for 2021-03-02 (serial date 44257) ---
$newRows = @(
    @("16-20",   7990,   0.556535742619981, 5261,   0.563728904366461, 2718,   0.545514756796355),
    @("21-30",   103730, 7.22521308910771,  63217,  6.77385480846504,  40383,  8.10504872101075),
    @("31-40",   136180, 9.48548653691977,  82958,  8.88915081703724,  52965,  10.6303123965118),
    @("41-50",   149071, 10.3833967068965,  92093,  9.86798821323332,  56617,  11.3632851308091),
    @("51-60",   165141, 11.5027370553199,  103092, 11.0465577283686,  61651,  12.3736314463768),
    @("61-70",   234761, 16.3520509978985,  179048, 19.185427270292,   55177,  11.0742706901223),
    @("71-80",   412161, 28.7086768728403,  275135, 29.4813822662738,  135587, 27.2129173398629),
    @("81+",     226525, 15.7783803625771,  132352, 14.1818376640771,  93133,  18.6922096558922),
    @("PENDING", 108,    0.00752263582014492, 94,  0.0100723278864184, 14,  0.00280986261778844)
)

$startRow = 673
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    # Copy date cell format (style index) from the row above it so the new
    # date cells share the same style as the rest of column B.
    $ws.Range("B672").Copy() | Out-Null
    $ws.Range("B$r").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = 44257
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = $row[2]
    $ws.Cells.Item($r, 5).Value = $row[3]
    $ws.Cells.Item($r, 6).Value = $row[4]
    $ws.Cells.Item($r, 7).Value = $row[5]
    $ws.Cells.Item($r, 8).Value = $row[6]
}

# --- Update dimension / defined name range ---
$wb.Names.Item("FINAL_AGEGROUPS").RefersTo = "='FINAL_AGEGROUPS'!`$A`$1:`$H`$681"
